# Skype meeting notes (05/15): strike through people who didn't actually
# attend / confirm, and add a line noting Ilia's apology.
#
# Target shape (see diff):
#   P1: "Attended by: Angel, Georgi, Ilia, Mikaeil"
#       - paragraph mark itself gets <w:strike/>
#       - "Angel", "Georgi", "Ilia", "Mikaeil" each individually struck
#         (the rest of the text / separators stay plain)
#   P2 (new): a struck-through leading tab, then
#       "Ilia informed about not being able to attend"
#       (the _GoBack bookmark that used to sit in the "Planned activities:"
#       paragraph now lives here)
#   P3: "Planned activities:" (now its own clean paragraph, no bookmark)
#
# We build the first two paragraphs from scratch as literal OOXML and
# splice them in with Range.InsertXML, then delete the two original
# paragraphs they replace. This sidesteps the COM Font-on-a-collapsed-
# range quirks (applying strike to a zero-length range can bleed into
# neighboring paragraphs) and lets us control exactly which runs get
# <w:rPr><w:strike/></w:rPr> and which stay untouched.

$d = $word.ActiveDocument

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$firstPara = $d.Paragraphs(1)
$secondPara = $d.Paragraphs(2)

$insertPoint = $d.Range($firstPara.Range.Start, $firstPara.Range.Start)

$replacementXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:strike/></w:rPr></w:pPr><w:r><w:t>Attended by:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>Angel</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>Georgi</w:t></w:r><w:r><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>Ilia</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:strike/></w:rPr><w:t>Mikaeil</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:rPr><w:strike/></w:rPr><w:tab/></w:r><w:r><w:t>Ilia informed about not being able to attend</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t>Planned activities:</w:t></w:r></w:p>
'@

$insertPoint.InsertXML($replacementXml)

# After the insert, the three new paragraphs sit before the two original
# ones (old P1 "Attended by: ..." and old P2 "Planned activities:" with
# the bookmark) — remove that original pair now that it's been replaced.
$oldFirstPara = $d.Paragraphs(4)
$oldSecondPara = $d.Paragraphs(5)
$deleteRange = $d.Range($oldFirstPara.Range.Start, $oldSecondPara.Range.End)
$deleteRange.Delete()
